$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 287, shifting existing rows 287..360 down to 288..361.
$ws.Rows("287:287").Insert()

# Populate the newly inserted row 287 with the new record's data.
$ws.Range("A287").Value = 10
$ws.Range("B287").Value = "Vega Modelo de Temuco"
$ws.Range("C287").Value = "La Araucanía"
$ws.Range("D287").Value = 44855
$ws.Range("E287").Value = 9
$ws.Range("F287").Value = 100114013
$ws.Range("G287").Value = "Zanahoria"
$ws.Range("H287").Value = "Sin especificar"
$ws.Range("I287").Value = "Primera"
$ws.Range("J287").Value = 355
$ws.Range("K287").Value = 22000
$ws.Range("L287").Value = 25000
$ws.Range("M287").Value = 24070
$ws.Range("N287").Value = "$/saco 20 kilos"
$ws.Range("O287").Value = "Región del Bíobío"
$ws.Range("P287").Value = 1204
$ws.Range("Q287").Value = 20
$ws.Range("R287").Value = "Hortaliza"
